$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04347946725289499
$ws.Range("C2").Value = -0.1939596374787607
$ws.Range("B3").Value = -0.02569761831736754
$ws.Range("C3").Value = -0.3516465272167756
$ws.Range("B4").Value = 0.09119588401837582
$ws.Range("C4").Value = 0.5163541689608068
$ws.Range("B5").Value = 0.1408039466271324
$ws.Range("C5").Value = -0.3485913201265082
$ws.Range("B6").Value = 0.1367151699348035
$ws.Range("C6").Value = 0.03909464021695645
$ws.Range("B7").Value = 0.4256506979483463
$ws.Range("C7").Value = 0.1326950546932967
$ws.Range("B8").Value = 0.3873826698491664
$ws.Range("C8").Value = -0.389886956290172
$ws.Range("B9").Value = 0.07945036306217389
$ws.Range("C9").Value = 0.2535395216609476
$ws.Range("B10").Value = 0.5530058106659657
$ws.Range("C10").Value = -0.01388457528129947
$ws.Range("B11").Value = 0.1573830891862086
$ws.Range("C11").Value = 0.06562196948424935
$ws.Range("B12").Value = -0.3639087361843189
$ws.Range("C12").Value = 0.006692202857746783
$ws.Range("B13").Value = -0.3533773860940383
$ws.Range("C13").Value = -0.2690566977756325
$ws.Range("B14").Value = -0.1404970623387202
$ws.Range("C14").Value = 0.2831391637069788
$ws.Range("B15").Value = -0.07353604510650566
$ws.Range("C15").Value = -0.209460663690553
$ws.Range("B16").Value = 0.002709126943057505
$ws.Range("C16").Value = -0.1121177152412177
$ws.Range("B17").Value = -0.003893296495922256
$ws.Range("C17").Value = 0.04011349219002069
